{"js": "// LJRL-127 Agregar sysproperty para mostrar o no la funcionalidad de Club Lojack\n//\n// 1) Insert an (empty) bookmark right before the \"Activa o desactiva el uso\n//    de proxy en la comunicacion con las camaras.\" run (inside the\n//    \"camera.proxy:\" paragraph).\n// 2) Collapse the \"clubLoJack.url:\" paragraph's four runs into two runs and\n//    append a brand-new \"clubLoJack.show:\" paragraph right after it.\n\nconst body = context.document.body;\n\n// --- Step 1: insert the new empty bookmark -------------------------------\nconst proxyResults = body.search(\n  \"Activa o desactiva el uso de proxy en la comunicacion con las camaras.\",\n  { matchCase: true }\n);\nproxyResults.load(\"items\");\nawait context.sync();\n\nif (proxyResults.items.length === 0) {\n  throw new Error(\"Could not find the camera.proxy description text\");\n}\n\nconst proxyTextStart = proxyResults.items[0].getRange(\"Start\");\nproxyTextStart.insertBookmark(\"__DdeLink__226_327176666\");\nawait context.sync();\n\n// --- Step 2: rewrite the clubLoJack.url paragraph & add clubLoJack.show --\nconst clubResults = body.search(\"clubLoJack.url:\", { matchCase: true });\nclubResults.load(\"items\");\nawait context.sync();\n\nif (clubResults.items.length === 0) {\n  throw new Error(\"Could not find the clubLoJack.url: paragraph\");\n}\n\nconst clubParagraph = clubResults.items[0].paragraphs.getFirst();\nconst clubWholeRange = clubParagraph.getRange(\"Whole\");\n\nconst ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\">\n<pkg:xmlData>\n<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\"><Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/></Relationships>\n</pkg:xmlData>\n</pkg:part>\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>\n<w:p><w:pPr><w:pStyle w:val=\"style26\"/></w:pPr>\n<w:r><w:rPr><w:b/></w:rPr><w:t>clubLoJack.url:</w:t></w:r>\n<w:r><w:rPr/><w:t xml:space=\"preserve\"> URL de redireccion para ver los beneficios del club lo jack</w:t></w:r>\n</w:p>\n<w:p><w:pPr><w:pStyle w:val=\"style26\"/></w:pPr>\n<w:r><w:rPr><w:b/></w:rPr><w:t>clubLoJack.</w:t></w:r>\n<w:r><w:rPr><w:b/></w:rPr><w:t>show</w:t></w:r>\n<w:r><w:rPr><w:b/></w:rPr><w:t>:</w:t></w:r>\n<w:r><w:rPr/><w:t xml:space=\"preserve\"> </w:t></w:r>\n<w:r><w:rPr/><w:t>Si tiene el valor true activa la funcionalidad de club lo jack</w:t></w:r>\n</w:p>\n</w:body></w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>`;\n\nclubWholeRange.insertOoxml(ooxml, \"Replace\");\nawait context.sync();\n", "ps1": "# LJRL-127 Agregar sysproperty para mostrar o no la funcionalidad de Club Lojack\n#\n# 1) Insert an (empty) bookmark right before the \"Activa o desactiva el uso\n#    de proxy en la comunicacion con las camaras.\" run (inside the\n#    \"camera.proxy:\" paragraph).\n# 2) Collapse the \"clubLoJack.url:\" paragraph's four runs into two runs and\n#    append a brand-new \"clubLoJack.show:\" paragraph right after it.\n\n$d = $word.ActiveDocument\n\n# --- Step 1: insert the new empty bookmark --------------------------------\n$proxyRng = $d.Content\n$found = $proxyRng.Find.Execute(\"Activa o desactiva el uso de proxy en la comunicacion con las camaras.\")\nif (-not $found) {\n    throw \"Could not find the camera.proxy description text\"\n}\n$proxyRng.Collapse(1)  # wdCollapseStart\n$d.Bookmarks.Add(\"__DdeLink__226_327176666\", $proxyRng)\n\n# --- Step 2: rewrite the clubLoJack.url paragraph & add clubLoJack.show ---\n$clubRng = $d.Content\n$found2 = $clubRng.Find.Execute(\"clubLoJack.url:\")\nif (-not $found2) {\n    throw \"Could not find the clubLoJack.url: paragraph\"\n}\n$clubRng.Expand(4)  # wdParagraph - grabs the whole paragraph incl. mark\n\n$xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\">\n<pkg:xmlData>\n<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\"><Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/></Relationships>\n</pkg:xmlData>\n</pkg:part>\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>\n<w:p><w:pPr><w:pStyle w:val=\"style26\"/></w:pPr>\n<w:r><w:rPr><w:b/></w:rPr><w:t>clubLoJack.url:</w:t></w:r>\n<w:r><w:rPr/><w:t xml:space=\"preserve\"> URL de redireccion para ver los beneficios del club lo jack</w:t></w:r>\n</w:p>\n<w:p><w:pPr><w:pStyle w:val=\"style26\"/></w:pPr>\n<w:r><w:rPr><w:b/></w:rPr><w:t>clubLoJack.</w:t></w:r>\n<w:r><w:rPr><w:b/></w:rPr><w:t>show</w:t></w:r>\n<w:r><w:rPr><w:b/></w:rPr><w:t>:</w:t></w:r>\n<w:r><w:rPr/><w:t xml:space=\"preserve\"> </w:t></w:r>\n<w:r><w:rPr/><w:t>Si tiene el valor true activa la funcionalidad de club lo jack</w:t></w:r>\n</w:p>\n</w:body></w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>'\n\n$clubRng.InsertXML($xml)\n"}
